$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.Value = "'244.63"
$c.Style = "Normal"
$c = $ws.Cells.Item(3, 4)
$c.Value = "'23.92"
$c.Style = "Normal"
$c = $ws.Cells.Item(4, 4)
$c.Value = "'5.205"
$c.Style = "Normal"
$c = $ws.Cells.Item(5, 4)
$c.Value = "'0.05740"
$c.Style = "Normal"
$c = $ws.Cells.Item(6, 4)
$c.Value = "'6.473"
$c.Style = "Normal"
$c = $ws.Cells.Item(7, 4)
$c.Value = "'3.170"
$c.Style = "Normal"
$c = $ws.Cells.Item(8, 4)
$c.Value = "'0.8132"
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 4)
$c.Value = "'0.8689"
$c.Style = "Normal"
$c = $ws.Cells.Item(10, 4)
$c.Value = "'0.1369"
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 4)
$c.Value = "'0.06940"
$c.Style = "Normal"
$c = $ws.Cells.Item(12, 4)
$c.Value = "'0.03163"
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 4)
$c.Value = "'0.02921"
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 4)
$c.Value = "'0.09335"
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 4)
$c.Value = "'3.820"
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 4)
$c.Value = "'0.001536"
$c.Style = "Normal"
$c = $ws.Cells.Item(17, 4)
$c.Value = "'0.04700"
$c.Style = "Normal"
$ws.Cells.Item(18, 2).Value = "TigerCash"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$c = $ws.Cells.Item(18, 4)
$c.Value = "'0.006151"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "17TigerCashTCH"
$ws.Cells.Item(19, 2).Value = "BitKan"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$c = $ws.Cells.Item(19, 4)
$c.Value = "'0.001241"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "18BitKanKAN"
$ws.Cells.Item(20, 2).Value = "HotbitToken"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$c = $ws.Cells.Item(20, 4)
$c.Value = "'0.004105"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "19HotbitTokenHTB"
$ws.Cells.Item(21, 2).Value = "NitroEx"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$c = $ws.Cells.Item(21, 4)
$c.Value = "'0.00008701"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "20NitroExNTX"
$ws.Cells.Item(22, 2).Value = "LEO"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$c = $ws.Cells.Item(22, 4)
$c.Value = "'3.558"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "21LEOLEO"
$ws.Cells.Item(23, 2).Value = "BTSEToken"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$c = $ws.Cells.Item(23, 4)
$c.Value = "'2.159"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "22BTSETokenBTSE"
$ws.Cells.Item(24, 2).Value = "One"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$c = $ws.Cells.Item(24, 4)
$c.Value = "'0.01015"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "23OneONEBestin24h"
$c = $ws.Cells.Item(25, 4)
$c.Value = "'0.3185"
$c.Style = "Normal"
$c = $ws.Cells.Item(27, 4)
$c.Value = "'0.0002329"
$c.Style = "Normal"
$c = $ws.Cells.Item(40, 4)
$c.Value = "'0.03714"
$c.Style = "Normal"
$ws.Cells.Item(41, 2).Value = "KickToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$c = $ws.Cells.Item(41, 4)
$c.Value = "'0.006271"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "40KickTokenKICK"
$ws.Cells.Item(42, 2).Value = "BKEXToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$c = $ws.Cells.Item(42, 4)
$c.Value = "'0.1051"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "41BKEXTokenBKK"
$ws.Cells.Item(43, 2).Value = "CEJI"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$c = $ws.Cells.Item(43, 4)
$c.Value = "'0.002245"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "42CEJICEJIWorstin24h"
$c = $ws.Cells.Item(44, 4)
$c.Value = "'0.007465"
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 4)
$c.Value = "'0.00005472"
$c.Style = "Normal"
$c = $ws.Cells.Item(47, 4)
$c.Value = "'0.4539"
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 4)
$c.Value = "'0.003317"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "47BOLOBOLO"
